$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B13").Value = "Code module quản lý tài khoản"
